$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Solution")

$data = @(
    ,@("M1", "DO", "M1", "M1", "M1", "M1", "M3", "A1", "A1", "DO", "M1", "M1", "M1", "M3", "A1", "A1", "DO", "M1", "M1", "M3", "M1", "M1", "DO", "M1", "M1", "M1", "M3", "M1")
    ,@("A1", "A1", "DO", "A2", "A2", "M2", "A1", "M1", "A1", "A1", "DO", "M1", "M3", "A1", "DO", "M1", "M2", "A1", "A1", "M3", "M1", "A2", "M1", "M1", "A1", "DO", "M1", "A1")
    ,@("DO", "A1", "M1", "M1", "M3", "A1", "M3", "DO", "M3", "A1", "M1", "A1", "A1", "M1", "M1", "M3", "A1", "DO", "M1", "A1", "A1", "M1", "M3", "DO", "M1", "A1", "M1", "M1")
    ,@("DO", "M2", "A2", "A2", "M2", "M1", "M1", "DO", "M3", "M2", "A1", "M3", "M2", "M1", "M1", "DO", "M1", "A1", "M1", "M1", "M1", "M1", "A1", "A1", "DO", "M2", "A2", "A2")
    ,@("DO", "M1", "A1", "A1", "M1", "M2", "M1", "M1", "M1", "DO", "M1", "M1", "M1", "A2", "M1", "M2", "A2", "DO", "M2", "M2", "A1", "DO", "M1", "M1", "A2", "M1", "M2", "A2")
    ,@("A1", "A1", "M3", "DO", "A1", "A1", "A1", "A1", "A1", "M1", "DO", "A1", "A1", "M3", "A1", "A1", "DO", "M1", "A1", "A1", "M3", "A1", "A1", "A1", "DO", "A1", "A1", "M3")
    ,@("A2", "A1", "M1", "A2", "A1", "M2", "DO", "A2", "A1", "M1", "A2", "A2", "M2", "DO", "M1", "A2", "M2", "A2", "DO", "M2", "M1", "A2", "A1", "DO", "M1", "A2", "M1", "A2")
    ,@("DO", "M3", "A2", "M3", "M1", "A2", "A1", "DO", "M2", "M3", "M3", "M1", "M1", "A2", "A2", "DO", "M1", "M1", "M2", "A2", "A2", "M2", "A2", "M1", "M2", "A1", "M2", "DO")
    ,@("M1", "M1", "A2", "M1", "DO", "M2", "M1", "A2", "M2", "A2", "A2", "M1", "A1", "DO", "DO", "M1", "A2", "M2", "A2", "M1", "A2", "A1", "M2", "A1", "A1", "M1", "A2", "DO")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowData = $data[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowData[$c]
    }
}
